# Disaggregation of commodity Copper
# Rename the commodity label "Copper ores and concentrates" -> "Copper"
# across all year sheets (2000-2100), and apply the tiny value corrections
# that resulted from the underlying recalculation on a subset of sheets.

$wb = $excel.ActiveWorkbook

$oldLabel = "Copper ores and concentrates"
$newLabel = "Copper"

# Small floating point corrections to D4 on specific year sheets (row 4 = Copper)
$valueFixes = @{
    "2021" = 60783.98114902512
    "2025" = 115754.3298054403
    "2032" = 245046.0554110847
    "2041" = 910571.299199763
    "2048" = 3510062.265742251
    "2054" = 3773548.082014818
    "2056" = 3131634.983683897
    "2074" = 3542295.797384474
    "2077" = 3176048.268197267
    "2081" = 2845498.530304906
    "2085" = 3026569.828354888
    "2091" = 3734582.91319855
}

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("C4")
    if ($cell.Value2 -eq $oldLabel) {
        $cell.Value2 = $newLabel
    }

    if ($valueFixes.ContainsKey($ws.Name)) {
        $ws.Range("D4").Value2 = $valueFixes[$ws.Name]
    }
}
